# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets to reflect the latest scrape (each value incremented).

$wb = $excel.ActiveWorkbook

# Sheet "展览" - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 31
$ws1.Range("F5").Value = 24
$ws1.Range("F10").Value = 38
$ws1.Range("F11").Value = 1821
$ws1.Range("F13").Value = 115
$ws1.Range("F15").Value = 272
$ws1.Range("F18").Value = 18
$ws1.Range("F22").Value = 777
$ws1.Range("F23").Value = 312
$ws1.Range("F25").Value = 247

# Sheet "全部类型" - row => new F value (rows offset by the extra entry at F8)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 31
$ws4.Range("F5").Value = 24
$ws4.Range("F11").Value = 38
$ws4.Range("F12").Value = 1821
$ws4.Range("F14").Value = 115
$ws4.Range("F16").Value = 272
$ws4.Range("F19").Value = 18
$ws4.Range("F23").Value = 777
$ws4.Range("F24").Value = 312
$ws4.Range("F26").Value = 247
